# Apply the "DaysOff" rota update:
#  - new rota date
#  - clear the old "one off" day assignments (name moved out of col A/extra cols)
#  - shift the kitchen rota names up one slot and append RAP at the bottom
#  - clear the now-unused bottom legend row entries
#  - narrow column A now that it no longer holds long one-off text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header date
$ws.Range("B1").Value = "25/03/2020"

# Row 3: no more "one off" worker listed, the placeholder row A3 becomes a blank line
$ws.Range("A3").Value = "___________"
$ws.Range("C3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("H3").Value = ""

# Rows 4-7: the one-off column A entries are removed
$ws.Range("A4").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("A5").Value = ""
$ws.Range("A6").Value = ""
$ws.Range("A7").Value = ""

# Kitchen rota column B: shift names up one slot, new bottom entry RAP
$ws.Range("B6").Value = "MahaDeva"
$ws.Range("B8").Value = "Supriti"
$ws.Range("B9").Value = "Anna"
$ws.Range("B10").Value = "Dganit"
$ws.Range("B11").Value = "Shakti"
$ws.Range("B12").Value = "Anuka"
$ws.Range("B13").Value = "Ben"
$ws.Range("B14").Value = "Mahi"
$ws.Range("B15").Value = "RAP"

# Bottom legend row (19-20): clear the now unused one-off labels
$ws.Range("C19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("F19").Value = ""
$ws.Range("H19").Value = ""
$ws.Range("H20").Value = ""

# Column A is narrower now that it no longer carries the long one-off names
$ws.Columns.Item(1).ColumnWidth = 12.25
